$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3403
$ws.Range("I3").Value = 3488
$ws.Range("I4").Value = 821
$ws.Range("I5").Value = 326
$ws.Range("I6").Value = 3939
$ws.Range("I7").Value = 11977

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I7").Value = 396
$ws.Range("I8").Value = 742
$ws.Range("I11").Value = 190
$ws.Range("I15").Value = 146
$ws.Range("I16").Value = 32
$ws.Range("I19").Value = 316
$ws.Range("I21").Value = 65
$ws.Range("I22").Value = 31
$ws.Range("I23").Value = 112
$ws.Range("I25").Value = 55
$ws.Range("I29").Value = 777
$ws.Range("I30").Value = 43
$ws.Range("I33").Value = 532
$ws.Range("I36").Value = 163
$ws.Range("I37").Value = 385
$ws.Range("I42").Value = 414
$ws.Range("I44").Value = 90
$ws.Range("I45").Value = 24
$ws.Range("I47").Value = 82
$ws.Range("I48").Value = 157
$ws.Range("I49").Value = 96
$ws.Range("I50").Value = 53
$ws.Range("I51").Value = 116
$ws.Range("I52").Value = 260
$ws.Range("I53").Value = 130
$ws.Range("I54").Value = 267
$ws.Range("I57").Value = 51
$ws.Range("I63").Value = 50
$ws.Range("I65").Value = 263
$ws.Range("I67").Value = 467
$ws.Range("I68").Value = 39
$ws.Range("I73").Value = 102
$ws.Range("I77").Value = 66
$ws.Range("I78").Value = 166
$ws.Range("I79").Value = 310
$ws.Range("I85").Value = 547
$ws.Range("I89").Value = 135
$ws.Range("I90").Value = 149
$ws.Range("I91").Value = 145
$ws.Range("I93").Value = 66
$ws.Range("I94").Value = 107
$ws.Range("I101").Value = 11977

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 223
$ws.Range("I6").Value = 135
$ws.Range("I7").Value = 547

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 72
$ws.Range("I3").Value = 92
$ws.Range("I5").Value = 10
$ws.Range("I7").Value = 260

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 190

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 234
$ws.Range("I7").Value = 742

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I3").Value = 36
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 120
$ws.Range("I6").Value = 102
$ws.Range("I7").Value = 396

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I3").Value = 29
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 124
$ws.Range("I3").Value = 118
$ws.Range("I6").Value = 105
$ws.Range("I7").Value = 385

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 112
$ws.Range("I3").Value = 167
$ws.Range("I6").Value = 155
$ws.Range("I7").Value = 467

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 71
$ws.Range("I7").Value = 263

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 127
$ws.Range("I3").Value = 184
$ws.Range("I4").Value = 27
$ws.Range("I5").Value = 18
$ws.Range("I6").Value = 176
$ws.Range("I7").Value = 532

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 267

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 236
$ws.Range("I3").Value = 270
$ws.Range("I5").Value = 30
$ws.Range("I6").Value = 206
$ws.Range("I7").Value = 777

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 124
$ws.Range("I3").Value = 89
$ws.Range("I6").Value = 87
$ws.Range("I7").Value = 316

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 90

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 21
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 157

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 113
$ws.Range("I3").Value = 142
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 414

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 38
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 112

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 49
$ws.Range("I7").Value = 145

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I5").Value = 1
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 90
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 310

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 49
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 107

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I2").Value = 15
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 102

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 23
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 31

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I3").Value = 25
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = 24

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 32
